$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Refresh the crypto price/volume columns (D = Price, E = Volume(1h)).
# Cells whose new Price text reads as a plain number (e.g. "22.59") get their
# NumberFormat forced to Text ("@") first so Excel keeps storing the literal
# string instead of silently converting it to a numeric value.
$ws.Range("D2").Value = "27.155.27"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.677.34"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.59"
$ws.Range("E8").Value = "  +5.03%  "
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "1.915.57"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "1.682.03"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("E15").Value = "  +4.13%  "
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "27.101.88"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "235.68"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.85"
$ws.Range("E19").Value = "  -4.05%  "
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.53"
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.61"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.31"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  +0.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0502"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").Value = "1.537.10"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.601"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("E40").Value = "  +3.84%  "
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.03"
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").Value = "1.822.74"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.61"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("E48").Value = "  +3.55%  "
$ws.Range("E49").Value = "  +6.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.20"
$ws.Range("E50").Value = "  +3.28%  "
$ws.Range("E51").Value = "  -0.46%  "
